$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group_Company")
$ws.Columns("F:F").Insert()
Write-Host "Used range after insert:" $ws.UsedRange.Address()
Write-Host "F1:" $ws.Range("F1").Value()
Write-Host "G1:" $ws.Range("G1").Value()
